# Add a new row of stock data (2020-03-12) to the bottom of the sheet,
# matching the existing table's layout (timestamp, date, id, name,
# open, high, low, close, vol).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 37

$ws.Cells.Item($newRow, 1).Value = 1583971200
# Leading apostrophe forces these to stay text (like the rest of the
# column), instead of being parsed as a date / number.
$ws.Cells.Item($newRow, 2).Value = "'2020-03-12"
$ws.Cells.Item($newRow, 3).Value = "'0217"
$ws.Cells.Item($newRow, 4).Value = "PWRWELL"
$ws.Cells.Item($newRow, 5).Value = 0.26
$ws.Cells.Item($newRow, 6).Value = 0.265
$ws.Cells.Item($newRow, 7).Value = 0.25
$ws.Cells.Item($newRow, 8).Value = 0.255
$ws.Cells.Item($newRow, 9).Value = 7304700
